# Update gh-pages output data: "想去人数" (want-to-go count) figures refreshed
# for two events that appear on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition): row 2 -> F2 5345 -> 5353, row 7 -> F7 318 -> 322
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 5353
$wsExhibition.Range("F7").Value = 322

# Sheet "全部类型" (All types): row 2 -> F2 5345 -> 5353, row 8 -> F8 318 -> 322
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5353
$wsAll.Range("F8").Value = 322
